$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the (human-readable) header row, shifting
# the existing data rows down by one. Excel carries the formatting from
# the row above into the newly inserted row automatically.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the machine-friendly / technical column
# identifiers that hierarchically relate to the header in row 1 (this is
# what allows two columns to be related to build SKOS hierarchies).
$ws.Range("A2").Value = "intervalo-renta"
$ws.Range("B2").Value = "municipio-codigo"
$ws.Range("C2").Value = "ano"
$ws.Range("D2").Value = "municipio-nombre"

# The former trailing row (which only had "mapping-ano.xlsx" in column C)
# is no longer present in the final layout; after the insert it now sits
# at row 6, so remove it entirely.
$ws.Rows.Item(6).Delete()

$wb.Save()
